# The workbook's "UITestCases" sheet has a Yes/No dropdown in column E
# (data validation list "Yes,No"). Rows 2-28 currently read "No" (several
# of them additionally carrying leftover red/green "highlight" cell
# formatting from earlier ad-hoc styling). The edit flips all of those
# rows to "Yes" and normalizes their formatting back to the plain style
# already used by most of the column (and by the untouched rows 29-34),
# clearing the stray highlight fills.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the dropdown value for rows 2 through 28 to "Yes".
$ws.Range("E2:E28").Value = "Yes"

# Normalize formatting across E2:E28 to match the plain style already
# used by E2/E3/E5/... (and rows 29-34): copy that cell's format and
# paste it (format only) across the whole block. This clears the
# leftover red/green highlight fills from E4, E6:E8, E9, E15, E18, E20:E23,
# E25:E28 without disturbing the values we just set.
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E2:E28").PasteSpecial(-4122) | Out-Null
